# Configuration.xlsx update — "Update latest version 2"
# (Latest Version 2 on AMASS Website as of 1 May 2023)
#
# Content changes on the "setting_parameters" sheet:
#  - A14 header text: "Generating DATAQC report" -> "Generating Supplementary report"
#  - C1 / C5 / C14 header text: "Requirement" -> "Requirements"
#  - Column A widened to fit the new, longer header text
#  - Active selection moved from D19 to K9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Requirement" column header (repeated on each of the 3 section-header
# rows) is now plural "Requirements".
$ws.Range("C1").Value = "Requirements"
$ws.Range("C5").Value = "Requirements"
$ws.Range("C14").Value = "Requirements"

# Section header for the last block used to read "Generating DATAQC report";
# the report it produces is now labelled "Generating Supplementary report".
$ws.Range("A14").Value = "Generating Supplementary report"

# Column A needs to be widened to fit the longer header strings.
$ws.Columns.Item(1).ColumnWidth = 28.7

# Update the saved selection/active cell.
[void]$ws.Range("K9").Select()
